# Update core switch configuration to include VLAN200, the guest network

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newVlanList = "1,10,20,30,40,50,60,99,100,200"
$taggedNote  = "Tagged in case of aggregating switch at alliance station"

# All six VLAN-list cells in column E (rows 10-15) now show the same,
# updated list of VLANs, including the new VLAN200 guest network.
$ws.Range("E10:E15").Value = $newVlanList

# The "Tagged in case of aggregating switch..." note that used to live in
# column H for rows 14 and 15 moves over to column I.
$ws.Range("I14").Value = $taggedNote
$ws.Range("I14").Font.Italic = $true
$ws.Range("H14").Clear()

$ws.Range("I15").Value = $taggedNote
$ws.Range("I15").Font.Italic = $true
$ws.Range("H15").Clear()

# Update the active selection to match the saved workbook state.
$ws.Range("E12").Select()
